$d = $word.ActiveDocument

# --- locate the paragraph that ends with "Un termómetro." -----------------
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*term*metro.*") {
        $anchorIndex = $i
    }
}
$anchor = $d.Paragraphs($anchorIndex)

# Paragraphs used as formatting donors already present in the document:
#   - a numbered / bold "question" paragraph (e.g. paragraph 1)
#   - a plain "answer" paragraph (e.g. paragraph 2)
$questionDonor = $d.Paragraphs(1)
$answerDonor = $d.Paragraphs(2)

# 1) Blank paragraph right after the "termómetro" answer --------------------
$anchor.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs($anchorIndex + 1)

# 2) New numbered/bold question paragraph -----------------------------------
$insertPoint = $d.Range($blank1.Range.End, $blank1.Range.End)
$insertPoint.FormattedText = $questionDonor.Range.FormattedText
$questionPara = $d.Paragraphs($anchorIndex + 2)
$questionPara.Range.Text = "¿Qué limitaciones significativas se tendrían si no aplicamos un enfoque de sistemas? (5%)"

# 3) First answer paragraph --------------------------------------------------
$insertPoint = $d.Range($questionPara.Range.End, $questionPara.Range.End)
$insertPoint.FormattedText = $answerDonor.Range.FormattedText
$answerPara1 = $d.Paragraphs($anchorIndex + 3)
$answerPara1.Range.Text = "Una visión reduccionista del universo, centrada en objetos aislados."

# 4) Second answer paragraph -------------------------------------------------
$insertPoint = $d.Range($answerPara1.Range.End, $answerPara1.Range.End)
$insertPoint.FormattedText = $answerDonor.Range.FormattedText
$answerPara2 = $d.Paragraphs($anchorIndex + 4)
$answerPara2.Range.Text = "Aislamiento disciplinario."

# 5) Trailing blank paragraph ------------------------------------------------
$answerPara2.Range.InsertParagraphAfter()

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
